# Remove the trailing " (...)" organizational-unit suffix from the
# "Reported By" names on the "New Entries" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("New Entries")

$updates = @{
    2 = "Kostovski, Simon"
    3 = "Øhrgaard, Maria Karmark"
    4 = "Tsakova, Konstantina"
    5 = "Piedade, Pedro"
    6 = "Polanská, Pavlína"
    7 = "Prasad K, Nagendra"
    8 = "Wessner, Anna"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 9).Value = $updates[$row]
}
